$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price/volume data (GitHub Actions scheduled refresh)

$ws.Range("D2").Value = "'66.374.84"
$ws.Range("E2").Value = "  +0.24%  "

$ws.Range("D3").Value = "'3.246.94"
$ws.Range("E3").Value = "  +2.47%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'609.05"
$ws.Range("E5").Value = "  +0.43%  "

$ws.Range("D6").Value = "'156.89"
$ws.Range("E6").Value = "  +1.69%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "'3.246.13"
$ws.Range("E8").Value = "  +2.54%  "

$ws.Range("E9").Value = "  -0.36%  "

$ws.Range("E10").Value = "  +1.88%  "

$ws.Range("D11").Value = "'5.81"
$ws.Range("E11").Value = "  +3.83%  "

$ws.Range("E12").Value = "  -4.14%  "

$ws.Range("D13").Value = "'0.0000270"
$ws.Range("E13").Value = "  +0.96%  "

$ws.Range("D14").Value = "'39.01"
$ws.Range("E14").Value = "  +1.69%  "

$ws.Range("D15").Value = "'3.777.75"
$ws.Range("E15").Value = "  +2.40%  "

$ws.Range("D16").Value = "'66.499.97"
$ws.Range("E16").Value = "  +0.45%  "

$ws.Range("D17").Value = "'7.45"
$ws.Range("E17").Value = "  +0.49%  "

$ws.Range("D18").Value = "'3.246.07"
$ws.Range("E18").Value = "  +2.37%  "

$ws.Range("D19").Value = "'0.113"
$ws.Range("E19").Value = "  +1.19%  "

$ws.Range("D20").Value = "'504.42"
$ws.Range("E20").Value = "  -1.20%  "

$ws.Range("D21").Value = "'15.44"
$ws.Range("E21").Value = "  +0.38%  "

$ws.Range("D22").Value = "'0.750"
$ws.Range("E22").Value = "  +2.98%  "

$ws.Range("D23").Value = "'8.08"
$ws.Range("E23").Value = "  +0.88%  "

$ws.Range("D24").Value = "'14.68"
$ws.Range("E24").Value = "  -0.42%  "

$ws.Range("D25").Value = "'87.16"
$ws.Range("E25").Value = "  +2.93%  "

$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.05%  "

$ws.Range("D27").Value = "'3.02"
$ws.Range("E27").Value = "  +0.39%  "

$ws.Range("D28").Value = "'9.13"
$ws.Range("E28").Value = "  -0.05%  "

$ws.Range("D29").Value = "'2.39"
$ws.Range("E29").Value = "  -0.17%  "

$ws.Range("E30").Value = "  +46.51%  "

$ws.Range("D31").Value = "'6.97"
$ws.Range("E31").Value = "  -2.31%  "

$ws.Range("E32").Value = "  -5.88%  "

$ws.Range("D33").Value = "'27.93"
$ws.Range("E33").Value = "  -0.10%  "

$ws.Range("E34").Value = "  -0.08%  "

$ws.Range("E35").Value = "  -3.37%  "

$ws.Range("D36").Value = "'6.44"
$ws.Range("E36").Value = "  -1.15%  "

$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").Value = "'55.53"
$ws.Range("E37").Value = "  +1.26%  "

$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").Value = "'3.33"
$ws.Range("E38").Value = "  +18.68%  "

$ws.Range("D39").Value = "'0.0₃0785"
$ws.Range("E39").Value = "  +14.99%  "

$ws.Range("D40").Value = "'494.19"
$ws.Range("E40").Value = "  -1.47%  "

$ws.Range("D41").Value = "'0.0422"
$ws.Range("E41").Value = "  +0.55%  "

$ws.Range("E42").Value = "  +0.30%  "

$ws.Range("D43").Value = "'8.82"
$ws.Range("E43").Value = "  +0.90%  "

$ws.Range("D44").Value = "'0.293"
$ws.Range("E44").Value = "  -1.23%  "

$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").Value = "'2.51"
$ws.Range("E45").Value = "  +3.11%  "

$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "'2.982.13"
$ws.Range("E46").Value = "  +5.54%  "

$ws.Range("D47").Value = "'28.74"
$ws.Range("E47").Value = "  +2.71%  "

$ws.Range("D48").Value = "'2.51"
$ws.Range("E48").Value = "  +5.79%  "

$ws.Range("E49").Value = "  +2.02%  "

$ws.Range("E50").Value = "  -0.01%  "

$ws.Range("D51").Value = "'121.28"
$ws.Range("E51").Value = "  -0.86%  "
